$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6151345372200012
$ws.Range("B1").Value = 2.556543827056885
$ws.Range("C1").Value = 6.301563262939453
$ws.Range("D1").Value = 1.727394700050354
$ws.Range("E1").Value = 1.578240036964417
